$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 (66ad14ad-... entry) - update handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-23 04:49:13"
$wsZhCn.Range("H4").Value = "2016-03-23 04:49:44"

# de-de sheet: row 4 (66ad14ad-... entry) - update handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-23 04:49:17"
$wsDeDe.Range("H4").Value = "2016-03-23 04:49:50"
